$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wb.Worksheets.Item(1).Name = "GNG_TO-1650996101607332"
$wb.Worksheets.Item(2).Name = "NB_TO-16509961035913308"
$wb.Worksheets.Item(3).Name = "RS_TO-16509961035913308"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509961036393301"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509961037033315"

# --- Sheet 1 (GNG_TO) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509961015753045.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961015913374.csv"
$ws1.Range("B4").Value = "go_stims-16509961015913374.csv"
$ws1.Range("B5").Value = "GNG_stims-1650996101607332.csv"

# --- Sheet 2 (NB_TO) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_0-16509961018873012.csv"
$ws2.Range("B3").Value = "ZB-match_0-16509961017352948.csv"
$ws2.Range("B4").Value = "TB-16509961034713292.csv"
$ws2.Range("B5").Value = "OB-1650996102479317.csv"
$ws2.Range("B6").Value = "ZB-match_5-16509961019753022.csv"
$ws2.Range("B7").Value = "TB-16509961035593297.csv"
$ws2.Range("B8").Value = "TB-16509961029912965.csv"
$ws2.Range("B9").Value = "OB-16509961026473417.csv"
$ws2.Range("B10").Value = "OB-16509961021593323.csv"

# --- Sheet 3 (RS_TO) --- no cell data changes

# --- Sheet 4 (TOL_TO) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509961036073291.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961035913308.csv"
$ws4.Range("B4").Value = "MM_stims-16509961036232955.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961036073291.csv"
$ws4.Range("B6").Value = "MM_stims-16509961036393301.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961036232955.csv"

# --- Sheet 5 (vSAT_TO) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16509961036873302.csv"
$ws5.Range("B3").Value = "SAT_stims-16509961036393301.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509961036713355.csv"
$ws5.Range("B5").Value = "SAT_stims-16509961036553314.csv"
